$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 51,4
$arr[0,0] = 1346434501199610112
$arr[0,1] = 'Pemprov DKI Jakarta melalui Dinas Pendidikan Provinsi DKI Jakarta tetap memberlakukan pembelajaran dari rumah untuk seluruh sekolah di Provinsi DKI Jakarta pada semester genap Tahun Ajaran (TA) 2020/2021.
#jagajakarta #disdikdki #belajardarirumah #jakartatanggapcorona https://t.co/EuWoq5sJ3Q'
$arr[0,2] = 'SolikhahID'
$arr[0,3] = 'Tue Jan 05 12:32:55 +0000 2021'
$arr[1,0] = 1346411269314597120
$arr[1,1] = 'Pembahasan Soal Usaha Part 3 (Usaha oleh Gaya Membentuk Sudut)
#Usaha #Energi #SoalUsaha #UsahaGayaMembentukSudut #PembahasanSoalUsaha #BelajarDaring #PembelajaranDaring #BelajarDariRumah https://t.co/BGdvULjzUD'
$arr[1,2] = 'FisikaCeria'
$arr[1,3] = 'Tue Jan 05 11:00:36 +0000 2021'
$arr[2,0] = 1346395677492802048
$arr[2,1] = '#GenPrestasi pernah main ke Candi Prambanan?
Candi Prambanan dibangun pada abad ke-9 dengan tinggi 47 meter.
#BelajarBarengIndiHomeStudy #IndiHomeStudyByIndiHome #BelajarLebihMudah #BelajarTanpaBatas #BelajarDariRumah #dirumahaja #IPSSMP #SejarahSMP #CandiPrambanan #TryOut https://t.co/a0DKb18pT7'
$arr[2,2] = 'indihome_study'
$arr[2,3] = 'Tue Jan 05 09:58:38 +0000 2021'
$arr[3,0] = 1346306295968591872
$arr[3,1] = 'Kunci Jawaban Tema 6 Kelas 5, Buku Tematik Tema Panas dan Perpindahannya, Contoh Sumber Energi Panas
#bukutematik #belajardarirumah 
https://t.co/DOI98eYObG'
$arr[3,2] = 'tribunkaltim'
$arr[3,3] = 'Tue Jan 05 04:03:28 +0000 2021'
$arr[4,0] = 1346284291374919936
$arr[4,1] = '🎧 Pandemi covid-19 memaksa aktivitas belajar-mengajar berubah di tahun 2020 kemarin. Di semester baru sekolah nanti besar kemungkinan akan online jg, masih ingatkah apa saja tantangannya saat mulai Pembelajaran Jarak Jauh? Share dengan reply twit ini ya.. 😉 👌 #BelajarDariRumah https://t.co/DAkfAWfVdl'
$arr[4,2] = 'AyoGerakBareng'
$arr[4,3] = 'Tue Jan 05 02:36:02 +0000 2021'
$arr[5,0] = 1346279756317146880
$arr[5,1] = 'Pembatalan pembelajaran tatap muka semester genap ini merupakan langkah yang dinilai oleh banyak orang tua dan masyarakat sebagai tindakan yang bijak bagi keamanan kesehatan para siswa.
#agentoleransi
#salamtoleransi
#BelajarDariRumah
#sekolah
#kemendikbud
https://t.co/WXcGuedtUp'
$arr[5,2] = 'YukToleran'
$arr[5,3] = 'Tue Jan 05 02:18:01 +0000 2021'
$arr[6,0] = 1346279500439559936
$arr[6,1] = 'SOAL TVRI KELAS 1 SD, EPISODE 2: ATURAN DI RUMAH, BELAJAR DARI RUMAH HARI INI, SELASA 5 JANUARI 2021.
#BelajardariRumah #BelajardariRumahTVRI #belajardirumah 
https://t.co/XJrLzbzWdO'
$arr[6,2] = 'tribunkaltim'
$arr[6,3] = 'Tue Jan 05 02:17:00 +0000 2021'
$arr[7,0] = 1346275627238756096
$arr[7,1] = 'SOAL TVRI KELAS 1 SD, EPISODE 2: ATURAN DI RUMAH, BELAJAR DARI RUMAH HARI INI, SELASA 5 JANUARI 2021.
#BelajardariRumah #BelajardariRumahTVRI #belajardirumah 
https://t.co/XJrLzbRx5m'
$arr[7,2] = 'tribunkaltim'
$arr[7,3] = 'Tue Jan 05 02:01:36 +0000 2021'
$arr[8,0] = 1346275226284293888
$arr[8,1] = 'Cara Membuat Web &amp; Landing Page Tanpa Coding Berbasis CMS Wordpress
https://t.co/NQ6GgMoJQD
#dirumahaja #pakaimasker #cucitangan #JagaJarak #WordPress #landingpage #umkm #ingatpesan #jagakesehatan #patuhiprotokolkesehatan #IndonesiaSehat #BelajarDariRumah #bersamalawancorona'
$arr[8,2] = 'PotekantropusX'
$arr[8,3] = 'Tue Jan 05 02:00:01 +0000 2021'
$arr[9,0] = 1346269678805549056
$arr[9,1] = 'Hai #SahabatPerpusdikbud , ayo kita simak jadwal belajar dari rumah hari ini.
Persiapkan diri kalian, jika mengalami kesulitan jangan sungkan ajak orang tua untuk ikut mendampingi belajar ya.
Semangat!
#merdekabelajar
#belajardarirumah
#bersamahadapikorona 
#Perpusdikbud https://t.co/UDafe2837H'
$arr[9,2] = 'perpusdikbud'
$arr[9,3] = 'Tue Jan 05 01:37:58 +0000 2021'
$arr[10,0] = 1346089395850419968
$arr[10,1] = 'Walaupun sebagian masih belajar dari rumah, tetap semangat yaa!
#BelajarBarengIndiHomeStudy #IndiHomeStudyByIndiHome #BelajarLebihMudah #BelajarTanpaBatas #BelajarDariRumah #dirumahaja #LiburSekolah'
$arr[10,2] = 'indihome_study'
$arr[10,3] = 'Mon Jan 04 13:41:35 +0000 2021'
$arr[11,0] = 1346088387938139904
$arr[11,1] = 'Pemerintah pusat mengizinkan pemerintah daerah untuk melaksanakan pembelajaran tatap muka (PTM).
#BelajarDariRumah #belajaronline #daring #KBM #UPI #IKAUPI #lampung #lampostco
https://t.co/rb8U9YjcY0'
$arr[11,2] = 'lampostco'
$arr[11,3] = 'Mon Jan 04 13:37:35 +0000 2021'
$arr[12,0] = 1346073614454465024
$arr[12,1] = 'kepingin insekyur tpi dee gaisok masak🤣
rebus air ae gosong
#Insecure #blitar #BelajarDariRumah'
$arr[12,2] = 'humanseblak'
$arr[12,3] = 'Mon Jan 04 12:38:53 +0000 2021'
$arr[13,0] = 1346042820201598976
$arr[13,1] = 'Dinas Pendidikan Kabupaten Pesawaran, memutuskan kegiatan belajar mengajar (KBM) tatap muka di Bumi Andan Jejama diundur sampai waktu yang belum ditentukan.
#BelajarDariRumah #belajaronline #daring #KBM #pesawaran #lampung #lampostco
https://t.co/SQduyCyQxr'
$arr[13,2] = 'lampostco'
$arr[13,3] = 'Mon Jan 04 10:36:31 +0000 2021'
$arr[14,0] = 1345966134432943872
$arr[14,1] = 'Soal dan Jawaban TVRI Kelas 6 SD Senin 4 Januari 2021 Belajar dari Rumah tentang Alat Transportasi 
https://t.co/n06qxDIy7j via @tribunkaltim #BelajarDariRumah #SoalJawabanTVRI #AlatTransportasi'
$arr[14,2] = 'tribunkaltim'
$arr[14,3] = 'Mon Jan 04 05:31:47 +0000 2021'
$arr[15,0] = 1345940669911552000
$arr[15,1] = 'Yok semangat yok udah masuk aja sekolahnya tapi tetep Daring. Semester terakhir hmm tinggal nunggu ujian :)
#sekolah #daring #coronavirus #BelajarDariRumah #msteams #smk'
$arr[15,2] = '_saddhuman_'
$arr[15,3] = 'Mon Jan 04 03:50:36 +0000 2021'
$arr[16,0] = 1345936822845678080
$arr[16,1] = 'Kemendikbud telah menyiapkan program Belajar Dari Rumah (BDR) untuk mendukung PJJ sebagai alternatif pembelajaran di masa pandemi. 
Yuk, lihat jadwal pembelajaran program BDR bersama TVRI, Minggu pertama di Tahun 2021 
#MerdekaBelajar 
#BersamaHadapiKorona 
#BelajardariRumah https://t.co/K6EiUneJDl'
$arr[16,2] = 'LPMPBALI'
$arr[16,3] = 'Mon Jan 04 03:35:19 +0000 2021'
$arr[17,0] = 1345931200376946944
$arr[17,1] = 'Kelas dr rumah semester genap telah dimulai pagi ini. Dimulai dgn kelas luar bersama ; olahraga pagi 
Selamat menikmati kebersamaan bersama anak2  tuk penguatan akhlak, logika dan entrepreneur serta leadershipnya. #BelajarDariRumah https://t.co/VAzZsKdCDh'
$arr[17,2] = 'klik_ayet'
$arr[17,3] = 'Mon Jan 04 03:12:59 +0000 2021'
$arr[18,0] = 1345927669880218112
$arr[18,1] = 'Selamat pagi, #SobatKom
Yuk, simak juga jadwal tayang #BelajardariRumah di @tvrinasional untuk hari ini, Senin, 4 Januari 2021
Semangat belajar di tahun yang baru, ya, #SobatKom
#MerdekaBelajar #BersamaHadapiKorona https://t.co/VBJStK4e0M'
$arr[18,2] = 'kominfongw'
$arr[18,3] = 'Mon Jan 04 02:58:57 +0000 2021'
$arr[19,0] = 1345916114509828096
$arr[19,1] = 'Bagi mamah muda jangan galau dalam merawat anak era internet begini panduannya #BelajarDariRumah
https://t.co/bAvEcy4Cee'
$arr[19,2] = 'rinariris99'
$arr[19,3] = 'Mon Jan 04 02:13:02 +0000 2021'
$arr[20,0] = 1345912874485521920
$arr[20,1] = 'Selamat pagi teman kece 👋
Semangat buat semester barunya. 
Jangan lupa sarapan dan tetap jaga kesehatan kalian ya 🤗
#kejarcita #semesterbaru #BelajarDariRumah https://t.co/1bXfcGuWH4'
$arr[20,2] = 'kejarcitaid'
$arr[20,3] = 'Mon Jan 04 02:00:09 +0000 2021'
$arr[21,0] = 1345912834639719936
$arr[21,1] = 'Cara Membuat Web &amp; Landing Page Tanpa Coding Berbasis CMS Wordpress
https://t.co/NQ6GgMGlfd
#dirumahaja #pakaimasker #cucitangan #JagaJarak #WordPress #landingpage #umkm #ingatpesan #jagakesehatan #patuhiprotokolkesehatan #IndonesiaSehat #BelajarDariRumah #bersamalawancorona'
$arr[21,2] = 'PotekantropusX'
$arr[21,3] = 'Mon Jan 04 02:00:00 +0000 2021'
$arr[22,0] = 1345906862340805120
$arr[22,1] = 'SOAL DAN JAWABAN TVRI KELAS 4-6 HARI INI 4 JANUARI 2021, TEMA 5 ALAT TRANSPORTASI DAN JENIS-JENISNYA
#BelajarDariRumah #TribunKaltim 
https://t.co/qoW98mDcDe via @tribunkaltim'
$arr[22,2] = 'tribunkaltim'
$arr[22,3] = 'Mon Jan 04 01:36:16 +0000 2021'
$arr[23,0] = 1345902716682341888
$arr[23,1] = 'Jadwal Program Belajar dari Rumah di TVRI Minggu Pertama Tahun 2021.
#unggahulang #belajardarirumah 
#kemdikbud #lpmpkepri #kawanlpmpkepri https://t.co/eagCzLyksP'
$arr[23,2] = 'lpmpkepri'
$arr[23,3] = 'Mon Jan 04 01:19:47 +0000 2021'
$arr[24,0] = 1345901891008433920
$arr[24,1] = 'Jadwal Program "BELAJAR DARI RUMAH"  Hari Senin, 4 Januari 2021 pukul 08.00-11.30 WIB
#SemangatBaruMarta
#CharlesZoraZeeTylerMarta
#BelajarDariRumah
#MediaPemersatuBangsa
#TVRI #TVRINASIONAL https://t.co/I6U9pdqTVi'
$arr[24,2] = 'TheMartaSaputra'
$arr[24,3] = 'Mon Jan 04 01:16:31 +0000 2021'
$arr[25,0] = 1345897607667666944
$arr[25,1] = 'Hari ini (4/1/21), awal hari aktivitas bersekolah. Sekolah dari rumah sebagai pilihan tepat utk menekan penyebaran virus dan pengutamaan keselamatan generasi penerus bangsa. #BelajarDariRumah'
$arr[25,2] = 'siprinar'
$arr[25,3] = 'Mon Jan 04 00:59:29 +0000 2021'
$arr[26,0] = 1345890911176196096
$arr[26,1] = 'Jadwal #belajardarirumah di #TVRI Senin 4 Januari plus link #LiveStreaming 
https://t.co/a2UDIOnkyP'
$arr[26,2] = 'davidhinsatama'
$arr[26,3] = 'Mon Jan 04 00:32:53 +0000 2021'
$arr[27,0] = 1345886774585774080
$arr[27,1] = 'Kembali lagi ke aktivitas dimana kita harus mengerjakan tugas segunung 🤦
Welcome Back to Online School 
#BelajarDariRumah'
$arr[27,2] = 'rvandkh'
$arr[27,3] = 'Mon Jan 04 00:16:27 +0000 2021'
$arr[28,0] = 1345868018014768896
$arr[28,1] = 'SOAL TVRI KELAS 1-3 SD HARI INI 4 JANUARI 2021, TEMA 5 DENAH RUMAHKU, CONTOH PERATURAN DI RUMAH
#BDR #BelajarDariRumah #TribunKaltim
https://t.co/bfuceVScqP via @tribunkaltim'
$arr[28,2] = 'tribunkaltim'
$arr[28,3] = 'Sun Jan 03 23:01:55 +0000 2021'
$arr[29,0] = 1345767171339611904
$arr[29,1] = 'Jadwal Program "BELAJAR DARI RUMAH"  Hari Senin, 4 Januari 2021 pukul 08.00-11.30 WIB
#BelajarDariRumah
#MediaPemersatuBangsa
#TVRI #TVRINASIONAL https://t.co/OnhsMxbHh6'
$arr[29,2] = 'TVRINasional'
$arr[29,3] = 'Sun Jan 03 16:21:11 +0000 2021'
$arr[30,0] = 1345725088046300928
$arr[30,1] = 'Tetap #BelajarDariRumah dulu https://t.co/qdVwGpNBiN'
$arr[30,2] = 'thebaim'
$arr[30,3] = 'Sun Jan 03 13:33:58 +0000 2021'
$arr[31,0] = 1345714757362281984
$arr[31,1] = '#GenPrestasi ada yang besok sudah masuk sekolah? Atau masih libur seminggu lagi?
#BelajarBarengIndiHomeStudy #IndiHomeStudyByIndiHome #BelajarLebihMudah #BelajarTanpaBatas #BelajarDariRumah #dirumahaja #LiburSekolah https://t.co/xWoL6eTzzF'
$arr[31,2] = 'indihome_study'
$arr[31,3] = 'Sun Jan 03 12:52:54 +0000 2021'
$arr[32,0] = 1345683617188438016
$arr[32,1] = 'Gak perlu bingung lagi kalau banyak quiz atau ujian, belajar aja di Kelas Pintar!
#bimbel #bimbelonline #BelajarDariRumah #pjj #kelaspintarindonesia https://t.co/6oW9lVsIIC'
$arr[32,2] = 'kelaspintar_id'
$arr[32,3] = 'Sun Jan 03 10:49:10 +0000 2021'
$arr[33,0] = 1345654144783142912
$arr[33,1] = 'JADWAL Belajar dari Rumah TVRI Senin 4 Januari 2021 Semester Genap untuk PAUD SD Kelas 1 2 3 4 5 6
#belajardarirumah #belajartvri #semestergenap #belajartvri #tvri #belajar #senin4januari2021
 https://t.co/kTzwVZcCg8'
$arr[33,2] = 'tribunpontianak'
$arr[33,3] = 'Sun Jan 03 08:52:03 +0000 2021'
$arr[34,0] = 1345629736572050944
$arr[34,1] = 'Yuk, simak juga jadwal tayang #BelajardariRumah di @TVRINasional untuk hari Senin, 4 Januari 2021! Semangat belajar di tahun yang baru, ya, #SahabatDikbud! 
#BersamaHadapiKorona https://t.co/puk0kW1XSo'
$arr[34,2] = 'Kemdikbud_RI'
$arr[34,3] = 'Sun Jan 03 07:15:04 +0000 2021'
$arr[35,0] = 1345628893030486016
$arr[35,1] = 'Selamat siang, #SahabatDikbud. Dapatkan Panduan Pembelajaran Program #BelajardariRumah di @TVRINasional minggu pertama tahun 2021 dengan mengunjungi https://t.co/cYwQMWU6DE.
#BersamaHadapiKorona https://t.co/lZj6vc9exy'
$arr[35,2] = 'Kemdikbud_RI'
$arr[35,3] = 'Sun Jan 03 07:11:43 +0000 2021'
$arr[36,0] = 1345625884712082944
$arr[36,1] = '#SahabatEdukasi sudah siap mengikuti proses pembelajaran kembali di semester baru ini? Yuk lihat panduan pembelajaran program Belajar Dari Rumah (BDR) di TVRI untuk minggu pertama di tahun 2021.
https://t.co/jyBhHcWiia
#BelajarDariRumah #serubelajarkebiasaanbaru #merdekabelajar'
$arr[36,2] = 'pusdatin_dikbud'
$arr[36,3] = 'Sun Jan 03 06:59:46 +0000 2021'
$arr[37,0] = 1345609056497517056
$arr[37,1] = 'Libur telah usai, kini saatnya #SahabatPerpusdikbud bersiap kembali sekolah.
Bagi #SahabatPerpusdikbud yg memiliki anak sekolah PAUD dan SD yuk catat jadwalnya serta persiapkan diri untuk belajar dari rumah di @TVRINasional ya
#merdekabelajar
#belajardarirumah
#bersamahapaikorona https://t.co/hlpV39SKP2'
$arr[37,2] = 'perpusdikbud'
$arr[37,3] = 'Sun Jan 03 05:52:53 +0000 2021'
$arr[38,0] = 1345603967221341952
$arr[38,1] = 'Pemprov DKI Jakarta melalui Dinas Pendidikan Provinsi DKI Jakarta tetap memberlakukan pembelajaran dari rumah utk seluruh sekolah di Provinsi DKI Jakarta pada semester genap Tahun Ajaran 2020/2021.
https://t.co/oMSDnRdgzC
#jagajakarta #disdikdki #belajardarirumah #hadapibersama https://t.co/B43mgwYovf'
$arr[38,2] = 'DKIJakarta'
$arr[38,3] = 'Sun Jan 03 05:32:40 +0000 2021'
$arr[39,0] = 1345595744833647104
$arr[39,1] = 'Cara Membuat Web &amp; Landing Page Tanpa Coding Berbasis CMS Wordpress
https://t.co/NQ6GgMGlfd
#dirumahaja #pakaimasker #cucitangan #JagaJarak #WordPress #landingpage #umkm #ingatpesan #jagakesehatan #patuhiprotokolkesehatan #IndonesiaSehat #BelajarDariRumah #bersamalawancorona'
$arr[39,2] = 'PotekantropusX'
$arr[39,3] = 'Sun Jan 03 05:00:00 +0000 2021'
$arr[40,0] = 1345533343027192064
$arr[40,1] = 'Materi Konsep Usaha - #BelajarDariRumah https://t.co/OLJvzQ9vuR https://t.co/M56xqG7Xt1'
$arr[40,2] = 'FisikaCeria'
$arr[40,3] = 'Sun Jan 03 00:52:02 +0000 2021'
$arr[41,0] = 1345531189482799104
$arr[41,1] = 'Pembahasan Soal Usaha Part 2 (Usaha Oleh Gaya Mendatar)
#konsepusaha #usahagayamendatar #soalusaha #pembelajarandaring #belajardarirumah https://t.co/I37By0Y9S8'
$arr[41,2] = 'FisikaCeria'
$arr[41,3] = 'Sun Jan 03 00:43:28 +0000 2021'
$arr[42,0] = 1345524140191621120
$arr[42,1] = 'Berbekal dari pengalaman pribadi serta menyadari pentingnya makanan sehat untuk keluarga sehingga terangkum resep&amp; menu dalam buku ini. Kehadiran buku ini diharapkan membantu orang tua maupun mereka yg ingin terjun bisnis kuliner.
#MerdekaBelajar
#BelajarDariRumah
#Eperpusdikbud https://t.co/UIYmiiwdRf'
$arr[42,2] = 'perpusdikbud'
$arr[42,3] = 'Sun Jan 03 00:15:28 +0000 2021'
$arr[43,0] = 1345516767439831040
$arr[43,1] = 'Materi Fisika Kelas X : Materi Konsep Usaha - #BelajarDariRumah https://t.co/0SJ9s54STu'
$arr[43,2] = 'FisikaCeria'
$arr[43,3] = 'Sat Jan 02 23:46:10 +0000 2021'
$arr[44,0] = 1345496765936455936
$arr[44,1] = 'Jadwal Baru Belajar dari Rumah TVRI untuk PAUD dan SD Kelas 1-6, Tiap Senin-Jumat Selama 30 Menit
https://t.co/kk7OLxdr7G #JadwalBaru #BelajarDariRumah #TVRI #BelajarOnline'
$arr[44,2] = 'tribunkaltim'
$arr[44,3] = 'Sat Jan 02 22:26:41 +0000 2021'
$arr[45,0] = 1345398859829108992
$arr[45,1] = '@ProfesorZubairi SMA Xaverius 3 awal semester Genap masih  #BelajarDariRumah :) #StaySafe semuanya! https://t.co/pJBUYMTIZz'
$arr[45,2] = 'xavegaplg_'
$arr[45,3] = 'Sat Jan 02 15:57:39 +0000 2021'
$arr[46,0] = 1345315816976183040
$arr[46,1] = 'perpetual adalah sistem pencatatan persediaan yang di lakukan secara langsung pada jumlah dan harga pokoknya.
@agtn_rara 
#akuntansi
#persediaan
#vyb 
#anakakuntansi 
#smkbisa
#belajardarirumah'
$arr[46,2] = 'agtn_rara'
$arr[46,3] = 'Sat Jan 02 10:27:40 +0000 2021'
$arr[47,0] = 1345303264565154048
$arr[47,1] = 'Dalam sebuah perenungan,dalam hening kita akan mendapat jawaban. Goresan malam karya Suara Senja dapat dinikmati dengan klik link ini https://t.co/APo7ZK2B2f
#puisiduarasa #renungan #hening #wellbeing #kesadaran #BelajarDariRumah #samasamabelajar #mindfulness #duarasaofficial https://t.co/3jEfJPz0L3'
$arr[47,2] = 'DuaRasaofficial'
$arr[47,3] = 'Sat Jan 02 09:37:47 +0000 2021'
$arr[48,0] = 1345271117406031872
$arr[48,1] = 'Dinas Pendidikan DKI Jakarta mengumumkan masih menerapkan pembelajaran dari rumah pada semester genap Tahun Ajaran (TA) 2020/2021.  #BelajardariRumah https://t.co/ef6TaKBnMe'
$arr[48,2] = 'jpnncom'
$arr[48,3] = 'Sat Jan 02 07:30:02 +0000 2021'
$arr[49,0] = 1345268030893338880
$arr[49,1] = 'Wali Kota Kediri Keluarkan SE Tentang Pelaksanaan Pembelajaran dari Rumah #Belajardarirumah #CegahPenyebaranCovid19 #WaliKotaKediriAbdullahAbuBakar https://t.co/ya29qPae0A https://t.co/GTvCoaxBXD'
$arr[49,2] = 'superradioid'
$arr[49,3] = 'Sat Jan 02 07:17:47 +0000 2021'
$arr[50,0] = 1345233359086128896
$arr[50,1] = 'Cara Membuat Web &amp; Landing Page Tanpa Coding Berbasis CMS Wordpress
https://t.co/NQ6GgMoJQD
#dirumahaja #pakaimasker #cucitangan #JagaJarak #WordPress #landingpage #umkm #ingatpesan #jagakesehatan #patuhiprotokolkesehatan #IndonesiaSehat #BelajarDariRumah #bersamalawancorona'
$arr[50,2] = 'PotekantropusX'
$arr[50,3] = 'Sat Jan 02 05:00:00 +0000 2021'

$ws.Range("A328:D378").Value = $arr

$ws.Range("L371").Select()

Write-Host "Added rows 328 to 378"
